$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell F1, styled like the other header cells (bold, centered,
# top-aligned, thin border - same look as B1:E1)
$headerCell = $ws.Cells.Item(1, 6)
$headerCell.Value = "time_taken"
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108  # xlCenter
$headerCell.VerticalAlignment = -4160    # xlTop
$headerCell.Borders.LineStyle = 1        # xlContinuous

# Timestamps for each data row (written as plain text, like the other
# string columns in this sheet)
$timestamps = @(
    "2021-10-05 10:52:42.609899",
    "2021-10-05 10:52:42.609911",
    "2021-10-05 10:52:42.609915",
    "2021-10-05 10:52:42.609919",
    "2021-10-05 10:52:42.609922",
    "2021-10-05 10:52:42.609925",
    "2021-10-05 10:52:42.609928",
    "2021-10-05 10:52:42.609931",
    "2021-10-05 10:52:42.609934"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.NumberFormat = "@"
    $cell.Value = $timestamps[$i]
}
